# Updated LR-pair edge-weight metrics with refreshed TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.113818
$ws.Range("H2").Value = 0.341454
$ws.Range("I2").Value = 0.0003230180320166274
$ws.Range("J2").Value = 0.0003230180320166274
$ws.Range("M2").Value = 0.4214143333333333
$ws.Range("N2").Value = 1.264243
$ws.Range("O2").Value = 0.02434128610922473
$ws.Range("P2").Value = 0.02434128610922473
$ws.Range("Q2").Value = 0.04796453659133333
$ws.Range("R2").Value = 0.431680829322
$ws.Range("S2").Value = 0.000007862674335755441
$ws.Range("T2").Value = 0.000007862674335755441

$ws.Range("G3").Value = 0.113818
$ws.Range("H3").Value = 0.341454
$ws.Range("I3").Value = 0.0003230180320166274
$ws.Range("J3").Value = 0.0003230180320166274
$ws.Range("N3").Value = 36.386704
$ws.Range("O3").Value = 0.7005766871049885
$ws.Range("P3").Value = 0.7005766871049887
$ws.Range("Q3").Value = 1.380487291957333
$ws.Range("R3").Value = 12.424385627616
$ws.Range("S3").Value = 0.0002262989027453819
$ws.Range("T3").Value = 0.000226298902745382

$ws.Range("G4").Value = 0.113818
$ws.Range("H4").Value = 0.341454
$ws.Range("I4").Value = 0.0003230180320166274
$ws.Range("J4").Value = 0.0003230180320166274
$ws.Range("M4").Value = 4.762423333333333
$ws.Range("N4").Value = 14.28727
$ws.Range("O4").Value = 0.2750820267857866
$ws.Range("P4").Value = 0.2750820267857866
$ws.Range("Q4").Value = 0.5420494989533333
$ws.Range("R4").Value = 4.87844549058
$ws.Range("S4").Value = 0.00008885645493548997
$ws.Range("T4").Value = 0.00008885645493548997

$ws.Range("I5").Value = 0.9904058666599795
$ws.Range("J5").Value = 0.9904058666599794
$ws.Range("M5").Value = 0.4214143333333333
$ws.Range("N5").Value = 1.264243
$ws.Range("O5").Value = 0.02434128610922473
$ws.Range("P5").Value = 0.02434128610922473
$ws.Range("Q5").Value = 147.0641070255747
$ws.Range("R5").Value = 1323.576963230172
$ws.Range("S5").Value = 0.02410775256462524
$ws.Range("T5").Value = 0.02410775256462524

$ws.Range("I6").Value = 0.9904058666599795
$ws.Range("J6").Value = 0.9904058666599794
$ws.Range("N6").Value = 36.386704
$ws.Range("O6").Value = 0.7005766871049885
$ws.Range("P6").Value = 0.7005766871049887
$ws.Range("Q6").Value = 4232.713276928491
$ws.Range("S6").Value = 0.6938552609539935
$ws.Range("T6").Value = 0.6938552609539935

$ws.Range("I7").Value = 0.9904058666599795
$ws.Range("J7").Value = 0.9904058666599794
$ws.Range("M7").Value = 4.762423333333333
$ws.Range("N7").Value = 14.28727
$ws.Range("O7").Value = 0.2750820267857866
$ws.Range("P7").Value = 0.2750820267857866
$ws.Range("Q7").Value = 1661.978436410787
$ws.Range("R7").Value = 14957.80592769708
$ws.Range("S7").Value = 0.2724428531413607
$ws.Range("T7").Value = 0.2724428531413607

$ws.Range("G8").Value = 3.266752
$ws.Range("H8").Value = 9.800256
$ws.Range("I8").Value = 0.009271115308003845
$ws.Range("J8").Value = 0.009271115308003843
$ws.Range("M8").Value = 0.4214143333333333
$ws.Range("N8").Value = 1.264243
$ws.Range("O8").Value = 0.02434128610922473
$ws.Range("P8").Value = 0.02434128610922473
$ws.Range("Q8").Value = 1.376656116245333
$ws.Range("R8").Value = 12.389905046208
$ws.Range("S8").Value = 0.0002256708702637348
$ws.Range("T8").Value = 0.0002256708702637347

$ws.Range("G9").Value = 3.266752
$ws.Range("H9").Value = 9.800256
$ws.Range("I9").Value = 0.009271115308003845
$ws.Range("J9").Value = 0.009271115308003843
$ws.Range("N9").Value = 36.386704
$ws.Range("O9").Value = 0.7005766871049885
$ws.Range("P9").Value = 0.7005766871049887
$ws.Range("Q9").Value = 39.62211268846933
$ws.Range("R9").Value = 356.599014196224
$ws.Range("S9").Value = 0.006495127248249679
$ws.Range("T9").Value = 0.006495127248249678

$ws.Range("G10").Value = 3.266752
$ws.Range("H10").Value = 9.800256
$ws.Range("I10").Value = 0.009271115308003845
$ws.Range("J10").Value = 0.009271115308003843
$ws.Range("M10").Value = 4.762423333333333
$ws.Range("N10").Value = 14.28727
$ws.Range("O10").Value = 0.2750820267857866
$ws.Range("P10").Value = 0.2750820267857866
$ws.Range("Q10").Value = 15.55765594901333
$ws.Range("R10").Value = 140.01890354112
$ws.Range("S10").Value = 0.00255031718949043
$ws.Range("T10").Value = 0.00255031718949043

